$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'325.01"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "'4.12%"
$ws.Range("E2").NumberFormat = "General"
$ws.Range("D3").Value = "'39.88"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "'6.78%"
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E4").Value = "'1.77%"
$ws.Range("E4").NumberFormat = "General"
$ws.Range("D5").Value = "'0.08099"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "'2.51%"
$ws.Range("E5").NumberFormat = "General"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.532"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "'2.53%"
$ws.Range("E6").NumberFormat = "General"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.623"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "'4.45%"
$ws.Range("E7").NumberFormat = "General"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.925"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "'0.85%"
$ws.Range("E8").NumberFormat = "General"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.958"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "'-1.40%"
$ws.Range("E9").NumberFormat = "General"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9349"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "'1.21%"
$ws.Range("E10").NumberFormat = "General"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1300"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "'14.32%"
$ws.Range("E11").NumberFormat = "General"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1959"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "'3.67%"
$ws.Range("E12").NumberFormat = "General"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09210"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "'1.76%"
$ws.Range("E13").NumberFormat = "General"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03414"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "'2.65%"
$ws.Range("E14").NumberFormat = "General"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09545"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "'-0.75%"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001395"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "'1.20%"
$ws.Range("E16").NumberFormat = "General"
$ws.Range("D17").Value = "'0.006517"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "'11.45%"
$ws.Range("E17").NumberFormat = "General"
$ws.Range("D18").Value = "'3.357"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3535"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "'2.52%"
$ws.Range("E19").NumberFormat = "General"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'6.696"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "'12.62%"
$ws.Range("E20").NumberFormat = "General"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1328"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "'3.16%"
$ws.Range("E21").NumberFormat = "General"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2313"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "'-10.75%"
$ws.Range("E22").NumberFormat = "General"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04438"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "'1.68%"
$ws.Range("E23").NumberFormat = "General"
$ws.Range("D24").Value = "'0.001222"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "'-1.12%"
$ws.Range("E24").NumberFormat = "General"
$ws.Range("D25").Value = "'0.004356"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "'-6.08%"
$ws.Range("E25").NumberFormat = "General"
$ws.Range("E26").Value = "'-5.24%"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("D27").Value = "'0.0003993"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "'-0.01%"
$ws.Range("E27").NumberFormat = "General"
$ws.Range("D39").Value = "'0.02454"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "'8.45%"
$ws.Range("E39").NumberFormat = "General"
$ws.Range("D40").Value = "'0.05231"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "'2.83%"
$ws.Range("E40").NumberFormat = "General"
$ws.Range("D41").Value = "'0.007702"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "'3.33%"
$ws.Range("E41").NumberFormat = "General"
$ws.Range("D42").Value = "'0.1431"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "'5.92%"
$ws.Range("E42").NumberFormat = "General"
$ws.Range("D43").Value = "'0.008677"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "'-3.86%"
$ws.Range("E43").NumberFormat = "General"
$ws.Range("E44").Value = "'1.08%"
$ws.Range("E44").NumberFormat = "General"
$ws.Range("D45").Value = "'0.008131"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "'-5.90%"
$ws.Range("E45").NumberFormat = "General"
$ws.Range("D46").Value = "'0.00006617"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "'-1.39%"
$ws.Range("E46").NumberFormat = "General"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("E47").NumberFormat = "General"
$ws.Range("D48").Value = "'0.002854"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "'-13.27%"
$ws.Range("E48").NumberFormat = "General"
$ws.Range("E49").Value = "'148.05%"
$ws.Range("E49").NumberFormat = "General"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").NumberFormat = "General"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "'-0.01%"
$ws.Range("E51").NumberFormat = "General"
